# Add team record columns (Wins, Losses, Ties) to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (AD1:AF1), matching the look of the existing headers.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting (bold font, border, centered/top alignment) from an
# existing header cell onto the new header cells.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the team record for every data row (rows 2-41).
$lastRow = 41
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 93  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 69  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
